# Weekly update: a new price-survey row for "Arveja Verde" (Femacal de La
# Calera, Coquimbo) is reported. It belongs chronologically right before the
# existing row 18, so insert a new row there (pushing the old rows 18-43 down
# to 19-44, including the dimension) and fill it with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(18).Insert()

$ws.Cells.Item(18, 1).Value = 3
$ws.Cells.Item(18, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(18, 3).Value = "Coquimbo"
$ws.Cells.Item(18, 4).Value = 44536
$ws.Cells.Item(18, 5).Value = 5
$ws.Cells.Item(18, 6).Value = 100112022
$ws.Cells.Item(18, 7).Value = "Arveja Verde"
$ws.Cells.Item(18, 8).Value = "Perfection"
$ws.Cells.Item(18, 9).Value = "Primera"
$ws.Cells.Item(18, 10).Value = 81
$ws.Cells.Item(18, 11).Value = 27000
$ws.Cells.Item(18, 12).Value = 28000
$ws.Cells.Item(18, 13).Value = 27556
$ws.Cells.Item(18, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(18, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(18, 16).Value = 1102
$ws.Cells.Item(18, 17).Value = 25
$ws.Cells.Item(18, 18).Value = "Hortaliza"
